$wb = $excel.ActiveWorkbook

$sheetNames = @(
    "SoBCaICbIC-urbanresidential",
    "SoBCaICbIC-ruralresidential",
    "SoBCaICbIC-commercial"
)

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # The "ISIC 05T06" column (header in C1) is being split into two
    # separate ISIC categories: "ISIC 05" (coal mining) and "ISIC 06"
    # (oil and gas extraction). Insert a new blank column before D so
    # everything from the old D column onward shifts one column right,
    # then label the two header cells.
    $ws.Columns("D").Insert()

    $ws.Range("C1").Value = "ISIC 05"
    $ws.Range("D1").Value = "ISIC 06"
}
